$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.218995928764343
$ws.Range("B1").Value = 2.233052492141724
$ws.Range("C1").Value = 2.986017942428589
$ws.Range("D1").Value = 3.469530344009399
$ws.Range("E1").Value = 1.672581553459167
